$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.042429616583995
$ws.Range("D2").Value = 1.049457249154486
$ws.Range("E2").Value = 1.050228828920262
$ws.Range("F2").Value = 1.060249651280139
$ws.Range("I2").Value = 1.039097680877022
$ws.Range("J2").Value = 1.04750548468758
$ws.Range("K2").Value = 1.052214254234562
$ws.Range("L2").Value = 1.052983687805741
$ws.Range("M2").Value = 1.062976943930869
$ws.Range("N2").Value = 1.005712725503983

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.043401473473572
$ws.Range("D3").Value = 1.050203754475487
$ws.Range("E3").Value = 1.051068603394769
$ws.Range("F3").Value = 1.061115303175654
$ws.Range("I3").Value = 1.039264375099796
$ws.Range("J3").Value = 1.048123665262534
$ws.Range("K3").Value = 1.052773179295484
$ws.Range("L3").Value = 1.053635796209765
$ws.Range("M3").Value = 1.063656854833151

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.044030973681405
$ws.Range("D4").Value = 1.050687164070496
$ws.Range("E4").Value = 1.051612867765385
$ws.Range("F4").Value = 1.061676206222317
$ws.Range("I4").Value = 1.039371074807582
$ws.Range("J4").Value = 1.048523693946497
$ws.Range("K4").Value = 1.053134545177883
$ws.Range("L4").Value = 1.054057979907647
$ws.Range("M4").Value = 1.064096926758423

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.044295768638591
$ws.Range("D5").Value = 1.050890476723763
$ws.Range("E5").Value = 1.051841884308697
$ws.Range("F5").Value = 1.06191219226649
$ws.Range("I5").Value = 1.039415652461198
$ws.Range("J5").Value = 1.04869187080648
$ws.Range("K5").Value = 1.053286391388146
$ws.Range("L5").Value = 1.054235518857298
$ws.Range("M5").Value = 1.064281961423095

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.044340237822431
$ws.Range("D6").Value = 1.050924618886334
$ws.Range("E6").Value = 1.051880349337011
$ws.Range("F6").Value = 1.061951826012878
$ws.Range("I6").Value = 1.039423120873681
$ws.Range("J6").Value = 1.048720108691822
$ws.Range("K6").Value = 1.053311882779311
$ws.Range("L6").Value = 1.054265331497563
$ws.Range("M6").Value = 1.06431303118716

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.04403451128444
$ws.Range("D7").Value = 1.050689880403128
$ws.Range("E7").Value = 1.051615927080967
$ws.Range("F7").Value = 1.061679358764614
$ws.Range("I7").Value = 1.039371671552988
$ws.Range("J7").Value = 1.048525941116403
$ws.Range("K7").Value = 1.053136574438746
$ws.Range("L7").Value = 1.054060351985298
$ws.Range("M7").Value = 1.064099399090579

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.042757926224945
$ws.Range("D8").Value = 1.049709456460889
$ws.Range("E8").Value = 1.05051245275892
$ws.Range("F8").Value = 1.060542042359917
$ws.Range("I8").Value = 1.039154256348442
$ws.Range("J8").Value = 1.047714395988748
$ws.Range("K8").Value = 1.052403206178902
$ws.Range("L8").Value = 1.053204023404193
$ws.Range("M8").Value = 1.063206696534087

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.040513398000349
$ws.Range("D9").Value = 1.047984729623433
$ws.Range("E9").Value = 1.048574749625076
$ws.Range("F9").Value = 1.058543898928308
$ws.Range("I9").Value = 1.038762264257143
$ws.Range("J9").Value = 1.046284584716024
$ws.Range("K9").Value = 1.051108701003354
$ws.Range("L9").Value = 1.051696845312641
$ws.Range("M9").Value = 1.061634647816357

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.039020454730539
$ws.Range("D10").Value = 1.046836953737093
$ws.Range("E10").Value = 1.0472875791803
$ws.Range("F10").Value = 1.057215897810546
$ws.Range("I10").Value = 1.038495000693024
$ws.Range("J10").Value = 1.045331598191676
$ws.Range("K10").Value = 1.05024427479114
$ws.Range("L10").Value = 1.050693327463663
$ws.Range("M10").Value = 1.060587368537109

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.038374815080541
$ws.Range("D11").Value = 1.046340456978481
$ws.Range("E11").Value = 1.046731336892778
$ws.Range("F11").Value = 1.056641849017959
$ws.Range("I11").Value = 1.038377871890806
$ws.Range("J11").Value = 1.044919010028748
$ws.Range("K11").Value = 1.049869645326843
$ws.Range("L11").Value = 1.050259108238989
$ws.Range("M11").Value = 1.060134079555237

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.03813511913626
$ws.Range("D12").Value = 1.046156112237947
$ws.Range("E12").Value = 1.046524892106737
$ws.Range("F12").Value = 1.056428771170416
$ws.Range("I12").Value = 1.038334154962514
$ws.Range("J12").Value = 1.044765766652612
$ws.Range("K12").Value = 1.049730443213382
$ws.Range("L12").Value = 1.050097867695635
$ws.Range("M12").Value = 1.059965737623938

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.0381865291315
$ws.Range("D13").Value = 1.046195651335829
$ws.Range("E13").Value = 1.046569167571929
$ws.Range("F13").Value = 1.056474470305264
$ws.Range("I13").Value = 1.038343541889631
$ws.Range("J13").Value = 1.044798637397144
$ws.Range("K13").Value = 1.049760304695792
$ws.Range("L13").Value = 1.050132452157661
$ws.Range("M13").Value = 1.06000184617191

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.038354999216644
$ws.Range("D14").Value = 1.046325217420456
$ws.Range("E14").Value = 1.046714268659358
$ws.Range("F14").Value = 1.05662423287936
$ws.Range("I14").Value = 1.038374262518198
$ws.Range("J14").Value = 1.044906342668657
$ws.Range("K14").Value = 1.049858139810506
$ws.Range("L14").Value = 1.050245779063451
$ws.Range("M14").Value = 1.060120163743929

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.038458815573265
$ws.Range("D15").Value = 1.046405057511067
$ws.Range("E15").Value = 1.046803692581793
$ws.Range("F15").Value = 1.05671652639123
$ws.Range("I15").Value = 1.038393162690621
$ws.Range("J15").Value = 1.044972704828634
$ws.Range("K15").Value = 1.049918412919758
$ws.Range("L15").Value = 1.050315609879881
$ws.Range("M15").Value = 1.060193067088632

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.039063320954634
$ws.Range("D16").Value = 1.04686991520542
$ws.Range("E16").Value = 1.047324518680401
$ws.Range("F16").Value = 1.057254016410272
$ws.Range("I16").Value = 1.038502744655688
$ws.Range("J16").Value = 1.045358981688417
$ws.Range("K16").Value = 1.050269130915621
$ws.Range("L16").Value = 1.050722151805551
$ws.Range("M16").Value = 1.06061745594117

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.039442730000337
$ws.Range("D17").Value = 1.047161642466965
$ws.Range("E17").Value = 1.047651517603972
$ws.Range("F17").Value = 1.057591434345611
$ws.Range("I17").Value = 1.038571107579598
$ws.Range("J17").Value = 1.045601300093061
$ws.Range("K17").Value = 1.050489040228158
$ws.Range("L17").Value = 1.050977248784326
$ws.Range("M17").Value = 1.060883715572234

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.039664111412381
$ws.Range("D18").Value = 1.047331849971951
$ws.Range("E18").Value = 1.047842357642732
$ws.Range("F18").Value = 1.057788339401595
$ws.Range("I18").Value = 1.038610847215302
$ws.Range("J18").Value = 1.045742646163314
$ws.Range("K18").Value = 1.050617277964092
$ws.Range("L18").Value = 1.051126072430407
$ws.Range("M18").Value = 1.061039038544481

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.03973961004405
$ws.Range("D19").Value = 1.047389894417232
$ws.Range("E19").Value = 1.047907447285956
$ws.Range("F19").Value = 1.05785549498445
$ws.Range("I19").Value = 1.038624374438098
$ws.Range("J19").Value = 1.045790842480589
$ws.Range("K19").Value = 1.050660998316721
$ws.Range("L19").Value = 1.051176822483623
$ws.Range("M19").Value = 1.061092002718771

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.039402014870104
$ws.Range("D20").Value = 1.047130337924484
$ws.Range("E20").Value = 1.047616422625923
$ws.Range("F20").Value = 1.057555222769186
$ws.Range("I20").Value = 1.038563786875735
$ws.Range("J20").Value = 1.045575301009541
$ws.Range("K20").Value = 1.05046544930534
$ws.Range("L20").Value = 1.05094987619235
$ws.Range("M20").Value = 1.060855146538547

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.038305385562102
$ws.Range("D21").Value = 1.046287061320813
$ws.Range("E21").Value = 1.046671535346282
$ws.Range("F21").Value = 1.056580127396632
$ws.Range("I21").Value = 1.038365221853509
$ws.Range("J21").Value = 1.044874625859779
$ws.Range("K21").Value = 1.049829331125199
$ws.Range("L21").Value = 1.050212405785317
$ws.Range("M21").Value = 1.060085321335817

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.037616605453966
$ws.Range("D22").Value = 1.045757301810711
$ws.Range("E22").Value = 1.046078422044168
$ws.Range("F22").Value = 1.055967911568892
$ws.Range("I22").Value = 1.038239160979643
$ws.Range("J22").Value = 1.044434143328695
$ws.Range("K22").Value = 1.049429100925165
$ws.Range("L22").Value = 1.049749006052155
$ws.Range("M22").Value = 1.059601474592385

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.037981672789205
$ws.Range("D23").Value = 1.046038094964528
$ws.Range("E23").Value = 1.0463927497406
$ws.Range("F23").Value = 1.056292376227829
$ws.Range("I23").Value = 1.038306103234731
$ws.Range("J23").Value = 1.044667645445168
$ws.Range("K23").Value = 1.04964129645536
$ws.Range("L23").Value = 1.049994636363206
$ws.Range("M23").Value = 1.059857954095896

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.039420412042297
$ws.Range("D24").Value = 1.04714448295181
$ws.Range("E24").Value = 1.047632280205095
$ws.Range("F24").Value = 1.057571584928855
$ws.Range("I24").Value = 1.038567095204805
$ws.Range("J24").Value = 1.045587048858304
$ws.Range("K24").Value = 1.05047610912579
$ws.Range("L24").Value = 1.05096224459617
$ws.Range("M24").Value = 1.060868055598917

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.041093066077208
$ws.Range("D25").Value = 1.048430259418182
$ws.Range("E25").Value = 1.049074882334041
$ws.Range("F25").Value = 1.059059751938384
$ws.Range("I25").Value = 1.038864652280607
$ws.Range("J25").Value = 1.046654190871923
$ws.Range("K25").Value = 1.051443617341201
$ws.Range("L25").Value = 1.052086268635498
$ws.Range("M25").Value = 1.062040932901625
